# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" worksheet (with fund-holding detail data) right after
# the "总计" summary sheet, and updates the "总计" sheet so its top data row
# now reflects the new 2022-Q4 aggregate, pushing the older quarters down by
# one row (2022-Q3, 2022-Q2, 2022-Q1).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet, positioned after "总计" (sheet 1), i.e.
#    immediately before the current "2022-Q3" sheet.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# The (still) 2nd sheet is now the old "2022-Q3" detail sheet - use it as a
# formatting template for the header row & the numbered index column so the
# new sheet's styling matches its siblings exactly.
$q3Template = $wb.Worksheets.Item(3)

# Header row formatting (bold / centered / bordered -> style used on B1:H1).
$q3Template.Range("B1:H1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null

# Index column (A) formatting for all 27 data rows (template only has 10
# rows, so copy it then paste across the full A2:A28 destination range).
$q3Template.Range("A2:A11").Copy() | Out-Null
$q4.Range("A2:A28").PasteSpecial(-4122) | Out-Null

# Header labels.
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Columns B..G hold text values (fund codes / names / numbers-as-text) in
# the source data, so force text format up front to stop values that look
# numeric (e.g. "14.62", "011466") from being silently coerced to numbers.
$q4.Range("B2:G28").NumberFormat = "@"

$q4Data = @(
    @("100060", "富国高新技术产业混合", "14.62", "93.20", "2.43", "0.3553", 10),
    @("011466", "兴业医疗保健混合A", "4.60", "88.16", "5.11", "0.2351", 7),
    @("007345", "富国科技创新灵活配置混合", "8.84", "94.61", "2.55", "0.2254", 10),
    @("519170", "浦银安盛增长动力灵活配置混合A", "6.85", "85.43", "2.22", "0.1521", 9),
    @("008619", "永赢医药健康股票C", "1.78", "89.12", "6.74", "0.1200", 6),
    @("163001", "长信医疗保健行业灵活配置混合（LOF）", "1.98", "93.22", "5.49", "0.1087", 3),
    @("011467", "兴业医疗保健混合C", "2.12", "88.16", "5.11", "0.1083", 7),
    @("008618", "永赢医药健康股票A", "0.58", "89.12", "6.74", "0.0391", 6),
    @("003284", "中邮医药健康灵活配置混合", "0.78", "91.51", "3.08", "0.0240", 9),
    @("003513", "中邮消费升级灵活配置混合", "0.53", "90.06", "4.47", "0.0237", 3),
    @("001243", "博时中证淘金大数据100指数I", "2.16", "91.19", "0.90", "0.0194", 8),
    @("007518", "东方阿尔法优选混合A", "0.94", "92.44", "1.88", "0.0177", 9),
    @("001242", "博时中证淘金大数据100指数A", "1.52", "91.19", "0.90", "0.0137", 8),
    @("001563", "华富健康文娱灵活配置混合", "0.33", "90.72", "3.38", "0.0112", 9),
    @("000649", "长城久鑫灵活配置混合A", "0.46", "90.33", "2.18", "0.0100", 8),
    @("007519", "东方阿尔法优选混合C", "0.41", "92.44", "1.88", "0.0077", 9),
    @("002068", "东方多策略灵活配置混合C", "0.26", "87.87", "2.78", "0.0072", 7),
    @("013166", "东兴宸祥量化混合A", "0.38", "93.88", "1.23", "0.0047", 3),
    @("009327", "东兴兴晟混合A", "0.38", "79.79", "1.08", "0.0041", 6),
    @("015655", "富荣医药健康混合A", "0.14", "92.48", "1.84", "0.0026", 2),
    @("013154", "长信医疗保健行业灵活配置混合(LOF)C", "0.02", "93.22", "5.49", "0.0011", 3),
    @("013167", "东兴宸祥量化混合C", "0.08", "93.88", "1.23", "0.0010", 3),
    @("009328", "东兴兴晟混合C", "0.08", "79.79", "1.08", "0.0009", 6),
    @("400023", "东方多策略灵活配置混合A", "0.03", "87.87", "2.78", "0.0008", 7),
    @("014003", "浦银安盛增长动力灵活配置混合C", "0.03", "85.43", "2.22", "0.0007", 9),
    @("015656", "富荣医药健康混合C", "0.01", "92.48", "1.84", "0.0002", 2),
    @("017461", "长城久鑫灵活配置混合C", "0.00", "90.33", "2.18", "0.0000", 8)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = ($r - 2)          # A: running index 0..26
    $q4.Cells.Item($r, 2).Value = $row[0]           # B: 基金代码
    $q4.Cells.Item($r, 3).Value = $row[1]           # C: 基金名称
    $q4.Cells.Item($r, 4).Value = $row[2]           # D: 基金规模
    $q4.Cells.Item($r, 5).Value = $row[3]           # E: 股票总仓位
    $q4.Cells.Item($r, 6).Value = $row[4]           # F: 仓位占比
    $q4.Cells.Item($r, 7).Value = $row[5]           # G: 持有市值(亿元)
    $q4.Cells.Item($r, 8).Value = $row[6]           # H: 仓位排名 (numeric)
    $r = $r + 1
}

# Last row's "持有市值" is exactly 0 -> stored as a genuine number (0), not
# text, in the source data - fix that one cell back up.
$q4.Range("G28").NumberFormat = "General"
$q4.Range("G28").Value = 0

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: the new 2022-Q4 totals become the
#    top data row, and the previously-top rows (2022-Q3 / 2022-Q2 / 2022-Q1)
#    shift down by one row each. The numbered index column (A) already holds
#    the correct sequential values (0,1,2) in rows 2-4, so only a new A5 is
#    required for the newly-appended row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 27
$total.Range("D2").Value = 1.49

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 10
$total.Range("D3").Value = 0.4

$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.14

# New row 5 - copy A4's style (bold/centered/bordered) onto A5 first.
$total.Range("A4").Copy() | Out-Null
$total.Range("A5").PasteSpecial(-4122) | Out-Null

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.07000000000000001

# Restore "总计" as the active/selected sheet (unchanged by the source diff).
$total.Activate()

Write-Output "2022-Q4 sheet added and 总计 sheet updated"
